# Generate Report for Handback
# Updates the timestamp strings that record when the handoff/handback
# xliff report generation occurred.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" for the first file
$wsOverview.Range("G2").Value = "2016-08-18 07:04:29"

# de-de!H2 - "Correspond Handoff Datetime" for the first file
# (shared the same original timestamp text as Overview!G2, so it moves in lockstep)
$wsDeDe.Range("H2").Value = "2016-08-18 07:04:29"

# zh-cn!H2 - "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-08-18 07:04:24"

# zh-cn!K2 - "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-08-18 07:04:41"

# de-de!K2 - "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-08-18 07:04:48"
